# pollster, hox & hermit: ability changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All Hat And No Cattle")

# Pollster (row 34) - reworded ability
$ws.Range("C34").Value = "When your exile is called for, visit the storyteller privately and choose a number. If exactly that many players vote, everyone who voted dies."

# Hermit (row 15) - new first-night-only-good-player-alive ability
$ws.Range("C15").Value = "The first time only one good player is alive at dawn, you wake and choose a player to come back to life."

# Lemming (row 22) - trimmed trailing whitespace from ability text
$ws.Range("C22").Value = "If a Lemming dies, so do all other Lemmings. When Lemmings die, a dead minion might regain their ability for one day and night."

# Hox (row 32) - reworded ability
$ws.Range("C32").Value = "Each night*, choose a player to die. If you choose yourself then (after your death) minions wake to choose a new evil hox and you become good. You don't learn bluffs or minions. Minions know you and learn bluffs. [+1 Outsider]"

# Update the selected/active cell to match the author's final cursor position
$ws.Range("C32").Select()
